$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to 79b83687-3609-4f18-b3c8-7a982e7f2368.md
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# zh-cn sheet: row 3 status, handback datetime, clear error detail
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-10-18 12:28:44"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).AutoFit() | Out-Null

# de-de sheet: row 3 status, handback datetime, clear error detail
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-10-18 12:29:01"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).AutoFit() | Out-Null
